$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (e.g. "511.55"), so they stay text
# exactly like the other price cells in the sheet (e.g. "60.579.55").
$textCells = @("D5", "D6", "D8", "D10", "D11", "D16", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D29", "D31", "D32", "D33", "D35", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D48")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '60.579.55'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').Value = '2.614.29'
$ws.Range('E3').Value = '  -0.41%  '
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').Value = '511.55'
$ws.Range('E5').Value = '  +0.47%  '
$ws.Range('D6').Value = '154.73'
$ws.Range('E6').Value = '  -2.52%  '
$ws.Range('E7').Value = '  +0.35%  '
$ws.Range('D8').Value = '0.587'
$ws.Range('E8').Value = '  -2.63%  '
$ws.Range('D9').Value = '2.627.14'
$ws.Range('E9').Value = '  -1.19%  '
$ws.Range('D10').Value = '6.74'
$ws.Range('E10').Value = '  +4.16%  '
$ws.Range('D11').Value = '0.104'
$ws.Range('E11').Value = '  -0.79%  '
$ws.Range('E12').Value = '  -0.04%  '
$ws.Range('E13').Value = '  +1.31%  '
$ws.Range('D14').Value = '3.071.94'
$ws.Range('E14').Value = '  -0.43%  '
$ws.Range('D15').Value = '60.519.03'
$ws.Range('E15').Value = '  -0.05%  '
$ws.Range('D16').Value = '21.60'
$ws.Range('E16').Value = '  -0.57%  '
$ws.Range('E17').Value = '  -0.38%  '
$ws.Range('D18').Value = '2.621.54'
$ws.Range('D19').Value = '4.76'
$ws.Range('E19').Value = '  -1.04%  '
$ws.Range('D20').Value = '352.48'
$ws.Range('E20').Value = '  +1.06%  '
$ws.Range('D21').Value = '10.59'
$ws.Range('E21').Value = '  +0.53%  '
$ws.Range('D22').Value = '6.18'
$ws.Range('E22').Value = '  -0.44%  '
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('D24').Value = '60.64'
$ws.Range('E24').Value = '  +0.40%  '
$ws.Range('D25').Value = '0.423'
$ws.Range('E25').Value = '  -0.40%  '
$ws.Range('D26').Value = '0.166'
$ws.Range('E26').Value = '  -0.48%  '
$ws.Range('E27').Value = '  +0.69%  '
$ws.Range('D28').Value = '0.0₃0842'
$ws.Range('E28').Value = '  -3.42%  '
$ws.Range('D29').Value = '7.35'
$ws.Range('E29').Value = '  -2.93%  '
$ws.Range('E30').Value = '  +0.21%  '
$ws.Range('D31').Value = '19.43'
$ws.Range('E31').Value = '  -0.76%  '
$ws.Range('B32').Value = 'Monero'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D32').Value = '151.45'
$ws.Range('E32').Value = '  -3.37%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Value = '1.57'
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('E34').Value = '  +0.46%  '
$ws.Range('D35').Value = '3.98'
$ws.Range('E35').Value = '  -2.10%  '
$ws.Range('E36').Value = '  -2.37%  '
$ws.Range('D37').Value = '0.893'
$ws.Range('E37').Value = '  +5.58%  '
$ws.Range('D38').Value = '1.49'
$ws.Range('E38').Value = '  -1.16%  '
$ws.Range('D39').Value = '0.845'
$ws.Range('E39').Value = '  -1.64%  '
$ws.Range('D40').Value = '36.33'
$ws.Range('E40').Value = '  +2.92%  '
$ws.Range('D41').Value = '3.76'
$ws.Range('E41').Value = '  -0.63%  '
$ws.Range('D42').Value = '291.26'
$ws.Range('E42').Value = '  -5.65%  '
$ws.Range('D43').Value = '0.626'
$ws.Range('E43').Value = '  -2.21%  '
$ws.Range('D44').Value = '0.101'
$ws.Range('E44').Value = '  -0.10%  '
$ws.Range('E45').Value = '  +0.63%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').Value = '0.0555'
$ws.Range('E46').Value = '  -4.06%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '19.87'
$ws.Range('E47').Value = '  +0.89%  '
$ws.Range('D48').Value = '4.91'
$ws.Range('E48').Value = '  -0.30%  '
$ws.Range('E49').Value = '  -0.67%  '
$ws.Range('E50').Value = '  +0.15%  '
$ws.Range('D51').Value = '1.998.60'
$ws.Range('E51').Value = '  -3.09%  '
